{"js": "// Update each two-digit multiplication prompt in the practice table to its\n// new operands (per the authoring diff). Every \"old\" expression is a unique\n// string in this document, so searching the body for each one and replacing\n// the match's text unambiguously targets the right cell without disturbing\n// any other run, paragraph, or formatting.\nconst pairs = [\n  [\"24\u00d782=\", \"72\u00d769=\"],\n  [\"27\u00d711=\", \"65\u00d742=\"],\n  [\"98\u00d765=\", \"70\u00d744=\"],\n  [\"65\u00d791=\", \"68\u00d766=\"],\n  [\"94\u00d725=\", \"61\u00d750=\"],\n  [\"75\u00d711=\", \"11\u00d757=\"],\n  [\"23\u00d728=\", \"78\u00d730=\"],\n  [\"38\u00d769=\", \"51\u00d718=\"],\n  [\"43\u00d777=\", \"35\u00d782=\"],\n  [\"92\u00d768=\", \"76\u00d737=\"],\n  [\"94\u00d718=\", \"18\u00d786=\"],\n  [\"39\u00d749=\", \"97\u00d721=\"],\n  [\"87\u00d746=\", \"52\u00d713=\"],\n  [\"36\u00d754=\", \"63\u00d719=\"],\n  [\"58\u00d797=\", \"34\u00d743=\"],\n  [\"47\u00d750=\", \"94\u00d772=\"],\n  [\"59\u00d739=\", \"67\u00d757=\"],\n  [\"96\u00d748=\", \"56\u00d758=\"],\n  [\"68\u00d774=\", \"29\u00d780=\"],\n  [\"43\u00d717=\", \"28\u00d732=\"],\n  [\"11\u00d778=\", \"50\u00d737=\"],\n  [\"25\u00d745=\", \"88\u00d750=\"],\n  [\"40\u00d775=\", \"40\u00d762=\"],\n  [\"66\u00d778=\", \"75\u00d770=\"],\n  [\"61\u00d728=\", \"80\u00d769=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication prompt in the practice table to its\n# new operands. Every old value is unique within the document, so a plain\n# Find/Replace (Replace:=wdReplaceAll, MatchCase, no wildcards) targets the\n# exact cell each time without touching anything else.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"24\u00d782=\", \"72\u00d769=\"),\n    @(\"27\u00d711=\", \"65\u00d742=\"),\n    @(\"98\u00d765=\", \"70\u00d744=\"),\n    @(\"65\u00d791=\", \"68\u00d766=\"),\n    @(\"94\u00d725=\", \"61\u00d750=\"),\n    @(\"75\u00d711=\", \"11\u00d757=\"),\n    @(\"23\u00d728=\", \"78\u00d730=\"),\n    @(\"38\u00d769=\", \"51\u00d718=\"),\n    @(\"43\u00d777=\", \"35\u00d782=\"),\n    @(\"92\u00d768=\", \"76\u00d737=\"),\n    @(\"94\u00d718=\", \"18\u00d786=\"),\n    @(\"39\u00d749=\", \"97\u00d721=\"),\n    @(\"87\u00d746=\", \"52\u00d713=\"),\n    @(\"36\u00d754=\", \"63\u00d719=\"),\n    @(\"58\u00d797=\", \"34\u00d743=\"),\n    @(\"47\u00d750=\", \"94\u00d772=\"),\n    @(\"59\u00d739=\", \"67\u00d757=\"),\n    @(\"96\u00d748=\", \"56\u00d758=\"),\n    @(\"68\u00d774=\", \"29\u00d780=\"),\n    @(\"43\u00d717=\", \"28\u00d732=\"),\n    @(\"11\u00d778=\", \"50\u00d737=\"),\n    @(\"25\u00d745=\", \"88\u00d750=\"),\n    @(\"40\u00d775=\", \"40\u00d762=\"),\n    @(\"66\u00d778=\", \"75\u00d770=\"),\n    @(\"61\u00d728=\", \"80\u00d769=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
